$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7167
$ws.Range("C2").Value = 250
$ws.Range("C3").Value = 151241
$ws.Range("C4").Value = 142891
$ws.Range("C5").Value = 8350
$ws.Range("C6").Value = 510
$ws.Range("C7").Value = 5.52
$ws.Range("C8").Value = 63.51
